# Auto update: 2025-12-03 03:05:24
# Refresh DECISION/국장_반도체_분석.xlsx with the latest pull of macro/ticker data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A ("날짜") is stored as plain text, e.g. "2025-12-03" -----------
# Excel auto-parses a bare "2025-12-03" assignment into a date serial, so
# force the cells to Text format first to keep the literal string value.
$ws.Range("A2:A7").NumberFormat = "@"
$ws.Range("A2:A7").Value = "2025-12-03"

# --- Row 2 : 058470.KS (SK hynix's old slot is replaced with DB HiTek's
#             ticker bundle in this refresh) ------------------------------
$ws.Cells.Item(2, 2).Value = "058470.KS,0P0000ASU1,98886"
$ws.Cells.Item(2, 3).Value = "058470.KS"
$ws.Cells.Item(2, 4).Value = 66600
$ws.Cells.Item(2, 5).Value = 66.90000000000001
$ws.Cells.Item(2, 6).Value = 18.09
$ws.Cells.Item(2, 7).Value = 60
$ws.Cells.Item(2, 8).Value = 50
$ws.Cells.Item(2, 9).Value = 63
$ws.Cells.Item(2, 10).Value = 70
$ws.Cells.Item(2, 11).Value = 62.8
$ws.Cells.Item(2, 14).Value = 65.32892478746797
$ws.Cells.Item(2, 15).Value = "🟢 상승 우위 (다소 완화)"

# --- Row 3 : SamsungElec ----------------------------------------------------
$ws.Cells.Item(3, 2).Value = "SamsungElec"
$ws.Cells.Item(3, 3).Value = "005930.KS"
$ws.Cells.Item(3, 4).Value = 103400
$ws.Cells.Item(3, 5).Value = 50.4
$ws.Cells.Item(3, 6).Value = 4.13
$ws.Cells.Item(3, 7).Value = 50
$ws.Cells.Item(3, 8).Value = 63
$ws.Cells.Item(3, 9).Value = 60
$ws.Cells.Item(3, 10).Value = 76
$ws.Cells.Item(3, 11).Value = 58.6
$ws.Cells.Item(3, 14).Value = 65.32892478746797
$ws.Cells.Item(3, 15).Value = "🟢 상승 우위 (다소 완화)"

# --- Row 4 : 403870.KS ------------------------------------------------------
$ws.Cells.Item(4, 2).Value = "403870.KS,0P0001PE9K,566428"
$ws.Cells.Item(4, 3).Value = "403870.KS"
$ws.Cells.Item(4, 4).Value = 31300
$ws.Cells.Item(4, 5).Value = 46.2
$ws.Cells.Item(4, 6).Value = 10.99
$ws.Cells.Item(4, 7).Value = 50
$ws.Cells.Item(4, 8).Value = 43
$ws.Cells.Item(4, 9).Value = 53
$ws.Cells.Item(4, 10).Value = 60
$ws.Cells.Item(4, 11).Value = 55.8
$ws.Cells.Item(4, 14).Value = 65.32892478746797
$ws.Cells.Item(4, 15).Value = "🟢 상승 우위 (다소 완화)"

# --- Row 5 : DB HiTek --------------------------------------------------------
$ws.Cells.Item(5, 2).Value = "DB HiTek"
$ws.Cells.Item(5, 3).Value = "000990.KS"
$ws.Cells.Item(5, 4).Value = 65000
$ws.Cells.Item(5, 5).Value = 35.3
$ws.Cells.Item(5, 6).Value = 5.01
$ws.Cells.Item(5, 7).Value = 30
$ws.Cells.Item(5, 8).Value = 53
$ws.Cells.Item(5, 9).Value = 60
$ws.Cells.Item(5, 10).Value = 60
$ws.Cells.Item(5, 11).Value = 52.6
$ws.Cells.Item(5, 14).Value = 65.32892478746797
$ws.Cells.Item(5, 15).Value = "🟢 상승 우위 (다소 완화)"

# --- Row 6 : SK hynix --------------------------------------------------------
$ws.Cells.Item(6, 2).Value = "SK hynix"
$ws.Cells.Item(6, 3).Value = "000660.KS"
$ws.Cells.Item(6, 4).Value = 558000
$ws.Cells.Item(6, 5).Value = 39.4
$ws.Cells.Item(6, 6).Value = 7.59
$ws.Cells.Item(6, 7).Value = 20
$ws.Cells.Item(6, 8).Value = 66
$ws.Cells.Item(6, 9).Value = 63
$ws.Cells.Item(6, 10).Value = 73
$ws.Cells.Item(6, 11).Value = 50.8
$ws.Cells.Item(6, 14).Value = 65.32892478746797
$ws.Cells.Item(6, 15).Value = "🟢 상승 우위 (다소 완화)"

# --- Row 7 : 240810.KS -------------------------------------------------------
$ws.Cells.Item(7, 2).Value = "240810.KS,0P00017YB3,330568"
$ws.Cells.Item(7, 3).Value = "240810.KS"
$ws.Cells.Item(7, 4).Value = 61500
$ws.Cells.Item(7, 5).Value = 28.6
$ws.Cells.Item(7, 6).Value = 3.71
$ws.Cells.Item(7, 7).Value = 20
$ws.Cells.Item(7, 8).Value = 63
$ws.Cells.Item(7, 9).Value = 60
$ws.Cells.Item(7, 10).Value = 66
$ws.Cells.Item(7, 11).Value = 49.6
$ws.Cells.Item(7, 14).Value = 65.32892478746797
$ws.Cells.Item(7, 15).Value = "🟢 상승 우위 (다소 완화)"
